# "fix bug client add album"
# The roadmap sheet ("Sheet1") gets a new "Client" / "Client service" work
# item (with its 3 sub-bullets) inserted right after the "Song manage" row,
# and the old mis-typed "Accout" / "Accout manage" pair is turned back into
# the "Album" task that had been dropped from the checklist.

$wb = $excel.ActiveWorkbook
$wsModel   = $wb.Worksheets.Item("model")
$wsService = $wb.Worksheets.Item("service")
$wsTodo    = $wb.Worksheets.Item("Sheet1")

# --- make room: push rows 7+ down by two rows (keeps row 6 in place for now) ---
$wsTodo.Rows("7:8").Insert()

# --- row 4: "Album" -> "Client" (new task), add the date range in column D ---
$wsTodo.Cells.Item(4, 2).Value = "Client"
$wsTodo.Cells.Item(4, 4).Value = "8/11 - 9/11"

# --- row 5: "Accout" -> "Client service", add the date range in column D ---
$wsTodo.Cells.Item(5, 2).Value = "Client service"
$wsTodo.Cells.Item(5, 4).Value = "8/11 - 9/11"

# --- row 6: used to be "Accout manage" (A=5); becomes the first sub-bullet,
#     no longer carries an index number or the "done" note ---
$wsTodo.Cells.Item(6, 1).ClearContents()
$wsTodo.Cells.Item(6, 3).ClearContents()
$wsTodo.Cells.Item(6, 2).Value = "+ Login / Register"

# --- rows 7-8: two brand new sub-bullets ---
$wsTodo.Cells.Item(7, 2).Value = "+ Singleton"
$wsTodo.Cells.Item(8, 2).Value = "+ IO DataClient"

# give the three "+ ..." bullets the same (quote-prefixed) style already used
# for this kind of heading elsewhere in the workbook (service!D17)
$wsService.Range("D17").Copy()
$wsTodo.Range("B6:B8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- row 9 (old row 7, only had an index number): restore the "Album" task ---
$wsTodo.Cells.Item(9, 2).Value = "Album"
$wsTodo.Cells.Item(9, 3).Value = "Làm xong 1 chức năng thì thử luôn"

# --- selection / active tab bookkeeping ---
# Sheet1 (todo list) was the active tab before; now "model" is, scrolled so
# row 7 is in view, while Sheet1 itself remembers a selection on B10.
$wsTodo.Range("B10").Select()
$wsModel.Activate()
